$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Resolving-Mac" cluster label is renamed to "Inflammatory-Mac" (the
# macrophage sub-population was reclassified with the new TPM pipeline).
# The target-cluster assignment for rows 4/5 and 8/9 is also swapped: row 4
# (and 8) now refers to the Inflammatory-Mac cluster, while row 5 (and 9)
# refers back to MuSCs.
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("D9").Value = "MuSCs"

# Updated TPM-derived metrics for all data rows (2-9).
$ws.Range("I2").Value = 0.1818801724491279
$ws.Range("J2").Value = 0.2500781318045117
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08108666666666665
$ws.Range("N2").Value = 0.24326
$ws.Range("O2").Value = 0.02056149724823249
$ws.Range("P2").Value = 0.02939174694341418
$ws.Range("Q2").Value = 0.002638614191111111
$ws.Range("R2").Value = 0.02374752772
$ws.Range("S2").Value = 0.003739728665320794
$ws.Range("T2").Value = 0.007350233166079984
$ws.Range("I3").Value = 0.1818801724491279
$ws.Range("J3").Value = 0.2500781318045117
$ws.Range("O3").Value = 0.02295278771378271
$ws.Range("P3").Value = 0.03280999043916413
$ws.Range("S3").Value = 0.004174656987571023
$ws.Range("T3").Value = 0.008205061113550054
$ws.Range("I4").Value = 0.1818801724491279
$ws.Range("J4").Value = 0.2500781318045117
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2176356666666667
$ws.Range("N4").Value = 0.652907
$ws.Range("O4").Value = 0.05518681856388939
$ws.Range("P4").Value = 0.07888710565478799
$ws.Range("Q4").Value = 0.007082009683777779
$ws.Range("R4").Value = 0.06373808715400001
$ws.Range("S4").Value = 0.01003738807731894
$ws.Range("T4").Value = 0.01972794000561451
$ws.Range("I5").Value = 0.1818801724491279
$ws.Range("J5").Value = 0.2500781318045117
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.5543775
$ws.Range("N5").Value = 7.108755
$ws.Range("O5").Value = 0.9012988964740954
$ws.Range("P5").Value = 0.8589111569626339
$ws.Range("Q5").Value = 0.115661813435
$ws.Range("R5").Value = 0.69397088061
$ws.Range("S5").Value = 0.1639283987189172
$ws.Range("T5").Value = 0.2147948975192672
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.146372
$ws.Range("H6").Value = 0.292744
$ws.Range("I6").Value = 0.8181198275508721
$ws.Range("J6").Value = 0.7499218681954883
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.08108666666666665
$ws.Range("N6").Value = 0.24326
$ws.Range("O6").Value = 0.02056149724823249
$ws.Range("P6").Value = 0.02939174694341418
$ws.Range("Q6").Value = 0.01186881757333333
$ws.Range("R6").Value = 0.07121290543999999
$ws.Range("S6").Value = 0.0168217685829117
$ws.Range("T6").Value = 0.0220415137773342
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.146372
$ws.Range("H7").Value = 0.292744
$ws.Range("I7").Value = 0.8181198275508721
$ws.Range("J7").Value = 0.7499218681954883
$ws.Range("O7").Value = 0.02295278771378271
$ws.Range("P7").Value = 0.03280999043916413
$ws.Range("Q7").Value = 0.013249154324
$ws.Range("R7").Value = 0.079494925944
$ws.Range("S7").Value = 0.01877813072621168
$ws.Range("T7").Value = 0.02460492932561407
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.146372
$ws.Range("H8").Value = 0.292744
$ws.Range("I8").Value = 0.8181198275508721
$ws.Range("J8").Value = 0.7499218681954883
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2176356666666667
$ws.Range("N8").Value = 0.652907
$ws.Range("O8").Value = 0.05518681856388939
$ws.Range("P8").Value = 0.07888710565478799
$ws.Range("Q8").Value = 0.03185576780133333
$ws.Range("R8").Value = 0.191134606808
$ws.Range("S8").Value = 0.04514943048657045
$ws.Range("T8").Value = 0.05915916564917348
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.146372
$ws.Range("H9").Value = 0.292744
$ws.Range("I9").Value = 0.8181198275508721
$ws.Range("J9").Value = 0.7499218681954883
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.5543775
$ws.Range("N9").Value = 7.108755
$ws.Range("O9").Value = 0.9012988964740954
$ws.Range("P9").Value = 0.8589111569626339
$ws.Range("Q9").Value = 0.5202613434300001
$ws.Range("R9").Value = 2.08104537372
$ws.Range("S9").Value = 0.7373704977551783
$ws.Range("T9").Value = 0.6441162594433667
